$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Lrfn3-Lrfn3)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lrfn3"
$ws.Cells.Item(2, 3).Value = "Lrfn3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.4264723333333333
$ws.Cells.Item(2, 8).Value = 1.279417
$ws.Cells.Item(2, 9).Value = 0.09324268104055088
$ws.Cells.Item(2, 10).Value = 0.1065665761548585
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.4264723333333333
$ws.Cells.Item(2, 14).Value = 1.279417
$ws.Cells.Item(2, 15).Value = 0.09324268104055088
$ws.Cells.Item(2, 16).Value = 0.1065665761548585
$ws.Cells.Item(2, 17).Value = 0.1818786510987778
$ws.Cells.Item(2, 18).Value = 1.636907859889
$ws.Cells.Item(2, 19).Value = 0.008694197567629906
$ws.Cells.Item(2, 20).Value = 0.01135643515336925

# Row 3: ECs -> FAPs (Lrfn3-Lrfn3)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lrfn3"
$ws.Cells.Item(3, 3).Value = "Lrfn3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.4264723333333333
$ws.Cells.Item(3, 8).Value = 1.279417
$ws.Cells.Item(3, 9).Value = 0.09324268104055088
$ws.Cells.Item(3, 10).Value = 0.1065665761548585
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.43175
$ws.Cells.Item(3, 14).Value = 7.295249999999999
$ws.Cells.Item(3, 15).Value = 0.531670806985587
$ws.Cells.Item(3, 16).Value = 0.6076438054940111
$ws.Cells.Item(3, 17).Value = 1.037074096583333
$ws.Cells.Item(3, 18).Value = 9.333666869249999
$ws.Cells.Item(3, 19).Value = 0.04957441147432938
$ws.Cells.Item(3, 20).Value = 0.06475451987320555

# Row 4: ECs -> sCs (Lrfn3-Lrfn3)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lrfn3"
$ws.Cells.Item(4, 3).Value = "Lrfn3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.4264723333333333
$ws.Cells.Item(4, 8).Value = 1.279417
$ws.Cells.Item(4, 9).Value = 0.09324268104055088
$ws.Cells.Item(4, 10).Value = 0.1065665761548585
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.7155665
$ws.Cells.Item(4, 14).Value = 3.431133
$ws.Cells.Item(4, 15).Value = 0.3750865119738621
$ws.Cells.Item(4, 16).Value = 0.2857896183511303
$ws.Cells.Item(4, 17).Value = 0.7316416482435
$ws.Cells.Item(4, 18).Value = 4.389849889461
$ws.Cells.Item(4, 19).Value = 0.03497407199859159
$ws.Cells.Item(4, 20).Value = 0.03045562112828366

# Row 5: FAPs -> ECs (Lrfn3-Lrfn3)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lrfn3"
$ws.Cells.Item(5, 3).Value = "Lrfn3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.43175
$ws.Cells.Item(5, 8).Value = 7.295249999999999
$ws.Cells.Item(5, 9).Value = 0.531670806985587
$ws.Cells.Item(5, 10).Value = 0.6076438054940111
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.4264723333333333
$ws.Cells.Item(5, 14).Value = 1.279417
$ws.Cells.Item(5, 15).Value = 0.09324268104055088
$ws.Cells.Item(5, 16).Value = 0.1065665761548585
$ws.Cells.Item(5, 17).Value = 1.037074096583333
$ws.Cells.Item(5, 18).Value = 9.333666869249999
$ws.Cells.Item(5, 19).Value = 0.04957441147432938
$ws.Cells.Item(5, 20).Value = 0.06475451987320555

# Row 6: FAPs -> FAPs (Lrfn3-Lrfn3)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lrfn3"
$ws.Cells.Item(6, 3).Value = "Lrfn3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.43175
$ws.Cells.Item(6, 8).Value = 7.295249999999999
$ws.Cells.Item(6, 9).Value = 0.531670806985587
$ws.Cells.Item(6, 10).Value = 0.6076438054940111
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.43175
$ws.Cells.Item(6, 14).Value = 7.295249999999999
$ws.Cells.Item(6, 15).Value = 0.531670806985587
$ws.Cells.Item(6, 16).Value = 0.6076438054940111
$ws.Cells.Item(6, 17).Value = 5.913408062499999
$ws.Cells.Item(6, 18).Value = 53.22067256249999
$ws.Cells.Item(6, 19).Value = 0.2826738470007053
$ws.Cells.Item(6, 20).Value = 0.3692309943552436

# Row 7: FAPs -> sCs (Lrfn3-Lrfn3)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lrfn3"
$ws.Cells.Item(7, 3).Value = "Lrfn3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.43175
$ws.Cells.Item(7, 8).Value = 7.295249999999999
$ws.Cells.Item(7, 9).Value = 0.531670806985587
$ws.Cells.Item(7, 10).Value = 0.6076438054940111
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.7155665
$ws.Cells.Item(7, 14).Value = 3.431133
$ws.Cells.Item(7, 15).Value = 0.3750865119738621
$ws.Cells.Item(7, 16).Value = 0.2857896183511303
$ws.Cells.Item(7, 17).Value = 4.171828836374999
$ws.Cells.Item(7, 18).Value = 25.03097301825
$ws.Cells.Item(7, 19).Value = 0.1994225485105523
$ws.Cells.Item(7, 20).Value = 0.1736582912655619

# Row 8: sCs -> ECs (Lrfn3-Lrfn3)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lrfn3"
$ws.Cells.Item(8, 3).Value = "Lrfn3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.7155665
$ws.Cells.Item(8, 8).Value = 3.431133
$ws.Cells.Item(8, 9).Value = 0.3750865119738621
$ws.Cells.Item(8, 10).Value = 0.2857896183511303
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.4264723333333333
$ws.Cells.Item(8, 14).Value = 1.279417
$ws.Cells.Item(8, 15).Value = 0.09324268104055088
$ws.Cells.Item(8, 16).Value = 0.1065665761548585
$ws.Cells.Item(8, 17).Value = 0.7316416482435
$ws.Cells.Item(8, 18).Value = 4.389849889461
$ws.Cells.Item(8, 19).Value = 0.03497407199859159
$ws.Cells.Item(8, 20).Value = 0.03045562112828366

# Row 9: sCs -> FAPs (Lrfn3-Lrfn3)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lrfn3"
$ws.Cells.Item(9, 3).Value = "Lrfn3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.7155665
$ws.Cells.Item(9, 8).Value = 3.431133
$ws.Cells.Item(9, 9).Value = 0.3750865119738621
$ws.Cells.Item(9, 10).Value = 0.2857896183511303
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.43175
$ws.Cells.Item(9, 14).Value = 7.295249999999999
$ws.Cells.Item(9, 15).Value = 0.531670806985587
$ws.Cells.Item(9, 16).Value = 0.6076438054940111
$ws.Cells.Item(9, 17).Value = 4.171828836374999
$ws.Cells.Item(9, 18).Value = 25.03097301825
$ws.Cells.Item(9, 19).Value = 0.1994225485105523
$ws.Cells.Item(9, 20).Value = 0.1736582912655619

# Row 10: sCs -> sCs (Lrfn3-Lrfn3)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lrfn3"
$ws.Cells.Item(10, 3).Value = "Lrfn3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.7155665
$ws.Cells.Item(10, 8).Value = 3.431133
$ws.Cells.Item(10, 9).Value = 0.3750865119738621
$ws.Cells.Item(10, 10).Value = 0.2857896183511303
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.7155665
$ws.Cells.Item(10, 14).Value = 3.431133
$ws.Cells.Item(10, 15).Value = 0.3750865119738621
$ws.Cells.Item(10, 16).Value = 0.2857896183511303
$ws.Cells.Item(10, 17).Value = 2.94316841592225
$ws.Cells.Item(10, 18).Value = 11.772673663689
$ws.Cells.Item(10, 19).Value = 0.1406898914647182
$ws.Cells.Item(10, 20).Value = 0.08167570595728468
